$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "produtividade de ELISSON MIGUEL GARCIA:5.9976"
$ws.Range("A3").Value = "produtividade de CARLOS ALBERTO BASILIO JUNIOR:5.6263"
$ws.Range("A4").Value = "produtividade de CLEMILSON SZNICER SOBRAL:5.2701"
$ws.Range("A5").Value = "produtividade de HEIDER DOS SANTOS NUNES:5.2509"
$ws.Range("A6").Value = "produtividade de EDUARDO FERREIRA LOURENCO:5.0893"
$ws.Range("A7").Value = "produtividade de TONY CRISTIAN DA SILVA JANDRE:4.8864"

$ws.Range("E7").Select()
